# Bind the active workbook/worksheet ourselves (the host-provided $wb can
# come through unbound in this runtime).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 was empty/unused (data started at row 2) - fill it in with the
# column titles used by the inventory parser, without shifting any of
# the existing data rows.
$ws.Range("A1").Value = "location"
$ws.Range("B1").Value = "shelf"
$ws.Range("C1").Value = "box"
$ws.Range("D1").Value = "item"
$ws.Range("E1").Value = "quantity"

# Move the active selection the same way the human edit left it.
$null = $ws.Range("F4").Select()
